# The document carries a "BTec_Logo-Orange" picture in both the primary
# and first-page headers, and a Pearson logo picture in both the primary
# and first-page footers. Each picture is an inline drawing; rename its
# drawing object (<wp:docPr name="...">, exposed on the COM object model
# as InlineShape.Name) the same way the source document was amended:
#   headers: image2.jpg -> image1.jpg
#   footers: image1.png -> image2.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($h = 1; $h -le 2; $h++) {
    $hdr = $sec.Headers.Item($h)
    $shapes = $hdr.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shapes.Item($i).Name = "image1.jpg"
    }
}

for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    $shapes = $ftr.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shapes.Item($i).Name = "image2.png"
    }
}
